$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cyclic shift of rows 2-4 (id, release_date):
#   new row2 = old row3
#   new row3 = old row4
#   new row4 = old row2
$oldA2 = $ws.Range("A2").Value2
$oldB2 = $ws.Range("B2").Value2
$oldA3 = $ws.Range("A3").Value2
$oldB3 = $ws.Range("B3").Value2
$oldA4 = $ws.Range("A4").Value2
$oldB4 = $ws.Range("B4").Value2

$ws.Range("A2").Value = $oldA3
$ws.Range("B2").Value = $oldB3

$ws.Range("A3").Value = $oldA4
$ws.Range("B3").Value = $oldB4

$ws.Range("A4").Value = $oldA2
$ws.Range("B4").Value = $oldB2
